# The commit swaps the deck's theme color palette from the "Integral"
# scheme over to the stock Office 2016+ "Office Theme" palette (the
# OOXML diff turns ppt/theme/theme1.xml's <a:clrScheme> from the green
# "Integral" values into the standard blue "Office" values).
#
# Re-point every slot of the active theme's 12-colour scheme
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) at the "Office
# Theme" palette used by PowerPoint's default Office theme.

$p = $ppt.ActivePresentation
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# Index order matches PowerPoint's msoThemeColorSchemeIndex layout:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink
$officeThemeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeThemeHex.Count; $i++) {
    $hex = $officeThemeHex[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($i + 1).RGB = $rgb
}
